# Apply updated audit numbers (Results sheet) and derived averages (Averages sheet).
$wb = $excel.ActiveWorkbook
$wsResults = $wb.Worksheets.Item("Results")

$resultsUpdates = @{
    2 = @{ "C"=50; "D"=12856.3315; "E"=31319.10899999998; "F"=12856.3315; "G"=5666949; "I"=1301930 }
    3 = @{ "D"=5875.24775; "E"=10476.48368800354; "F"=6733.031759310616; "G"=1537104; "I"=647102 }
    4 = @{ "C"=69; "D"=2406.10756; "E"=9846.369949999997; "F"=4477.799266667214; "G"=1183927; "I"=758093 }
    5 = @{ "C"=67; "D"=2536.668; "E"=10578.80974; "F"=4992.058988892054; "G"=1185604; "I"=650597 }
    6 = @{ "C"=62; "D"=4745.114; "E"=6105.682; "F"=5671.082166716559; "G"=1432910; "I"=770319 }
    7 = @{ "D"=6319.682999999999; "E"=13982.259; "F"=7257.048709542826; "G"=3220404; "I"=676829 }
    8 = @{ "C"=59; "D"=5837.446; "E"=8438.512999999999; "F"=5903.616913760536; "G"=1111324; "I"=432208 }
    9 = @{ "C"=64; "D"=4632.223999999999; "E"=5789.358999999999; "F"=4877.75103212671; "G"=1364378; "I"=770338 }
    10 = @{ "C"=63; "D"=4759.752; "E"=6078.044; "F"=5501.910514819657; "G"=1516659; "I"=770306 }
    11 = @{ "D"=3229.920000000001; "E"=10672.751; "F"=5713.540492995829; "G"=3112740; "I"=464810 }
    12 = @{ "C"=70; "D"=2820.24266; "E"=6461.98022; "F"=4208.519949127516; "G"=928329; "I"=248099 }
    13 = @{ "C"=55; "D"=9302.185; "E"=16804.97022; "F"=9302.185; "G"=2556234; "I"=830748 }
    14 = @{ "C"=75; "D"=2475.844000000001; "E"=6634.818; "F"=2475.844000000001; "G"=1239153; "I"=273758 }
    15 = @{ "C"=65; "D"=4557.028000000001; "E"=5759.991; "F"=5160.984846839408; "G"=1498904; "I"=770322 }
    16 = @{ "C"=63; "D"=4839.279; "E"=6519.782999999999; "F"=5347.843346807336; "G"=1708494; "I"=1109827 }
    17 = @{ "C"=58; "D"=6945.37916; "E"=16262.5585; "F"=6945.37916; "G"=6238283; "I"=757169 }
    18 = @{ "D"=6169.28; "E"=8051.973999999999; "F"=6169.28; "G"=2485573; "I"=1236728 }
    19 = @{ "C"=50; "D"=14623.8915; "E"=19100.0175; "F"=14623.8915; "G"=3158830; "I"=1624491 }
    20 = @{ "C"=56; "D"=8464.333999999999; "E"=23901.5495; "F"=8464.333999999999; "G"=6750718; "I"=602446 }
    21 = @{ "C"=56; "D"=8840.423; "E"=10867.421; "F"=8840.423; "G"=1595016; "I"=872009 }
    22 = @{ "C"=72; "D"=2771.625; "E"=6564.225; "F"=3032.918283528135; "G"=868685; "I"=386492 }
}

foreach ($row in $resultsUpdates.Keys) {
    $cols = $resultsUpdates[$row]
    foreach ($col in $cols.Keys) {
        $wsResults.Range("$col$row").Value = $cols[$col]
    }
}

$wsAverages = $wb.Worksheets.Item("Averages")
$averagesUpdates = @{
    2 = 61.62
    3 = 5952.76
    4 = 11438.89
    5 = 6597.89
    7 = 759743.86
    8 = 2.287
}

foreach ($row in $averagesUpdates.Keys) {
    $wsAverages.Range("B$row").Value = $averagesUpdates[$row]
}

Write-Output "Done updating Results and Averages sheets."